$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared-formula results in G1:G19 with static values
# (the formula (0.057*0.068) is being replaced by a fixed literal),
# and strip the centered-alignment style that these cells used to carry.
$rng = $ws.Range("G1:G19")
$rng.Value = 0.063455
$rng.ClearFormats()

# Move the active selection to I9, matching the latest user interaction.
$ws.Range("I9").Select() | Out-Null
